$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The PFAS selection criteria changed: "6:2 FTSA" and "PFPeA" are no longer
# included in the BMF diet computation, so their rows are removed from the
# sheet. All other rows stay the same; the sumPFAS row is recomputed.

# Delete the "6:2 FTSA" row (row 2).
$ws.Rows(2).Delete()

# After the first deletion, "PFPeA" (originally row 11) is now row 10.
$ws.Rows(10).Delete()

# Recompute the sumPFAS row (now row 13) with the new totals.
$ws.Range("B13").Value = 14.33
$ws.Range("C13").Value = 17.19
$ws.Range("D13").Value = 25.8
$ws.Range("E13").Value = 10.95
$ws.Range("F13").Value = 15.51
$ws.Range("G13").Value = 18.97
$ws.Range("H13").Value = 0.76
$ws.Range("I13").Value = 1.11
$ws.Range("J13").Value = 2.36
